$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows above the old row 146 (which holds the 44194 week's
# "Región de O'Higgins" data). This shifts the old rows 146-149 down to
# 150-153 unchanged, and leaves fresh blank rows at 146-149 for the new
# 44568 week's "Región del Maule" data.
$ws.Range("A146:A149").EntireRow.Insert()

# --- Row 146: Melón / Calameño / Extra ---
$ws.Cells.Item(146, 1).Value = 7
$ws.Cells.Item(146, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(146, 3).Value = "Ñuble"
$ws.Cells.Item(146, 4).Value = 44568
$ws.Cells.Item(146, 5).Value = 16
$ws.Cells.Item(146, 6).Value = 100112027
$ws.Cells.Item(146, 7).Value = "Melón"
$ws.Cells.Item(146, 8).Value = "Calameño"
$ws.Cells.Item(146, 9).Value = "Extra"
$ws.Cells.Item(146, 10).Value = 2000
$ws.Cells.Item(146, 11).Value = 800
$ws.Cells.Item(146, 12).Value = 900
$ws.Cells.Item(146, 13).Value = 850
$ws.Cells.Item(146, 14).Value = "$/unidad"
$ws.Cells.Item(146, 15).Value = "Región del Maule"
$ws.Cells.Item(146, 16).Value = 850
$ws.Cells.Item(146, 17).Value = 1
$ws.Cells.Item(146, 18).Value = "Hortaliza"

# --- Row 147: Melón / Calameño / Primera ---
$ws.Cells.Item(147, 1).Value = 7
$ws.Cells.Item(147, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(147, 3).Value = "Ñuble"
$ws.Cells.Item(147, 4).Value = 44568
$ws.Cells.Item(147, 5).Value = 16
$ws.Cells.Item(147, 6).Value = 100112027
$ws.Cells.Item(147, 7).Value = "Melón"
$ws.Cells.Item(147, 8).Value = "Calameño"
$ws.Cells.Item(147, 9).Value = "Primera"
$ws.Cells.Item(147, 10).Value = 5000
$ws.Cells.Item(147, 11).Value = 600
$ws.Cells.Item(147, 12).Value = 700
$ws.Cells.Item(147, 13).Value = 650
$ws.Cells.Item(147, 14).Value = "$/unidad"
$ws.Cells.Item(147, 15).Value = "Región del Maule"
$ws.Cells.Item(147, 16).Value = 650
$ws.Cells.Item(147, 17).Value = 1
$ws.Cells.Item(147, 18).Value = "Hortaliza"

# --- Row 148: Melón / Tuna / Extra ---
$ws.Cells.Item(148, 1).Value = 7
$ws.Cells.Item(148, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(148, 3).Value = "Ñuble"
$ws.Cells.Item(148, 4).Value = 44568
$ws.Cells.Item(148, 5).Value = 16
$ws.Cells.Item(148, 6).Value = 100112027
$ws.Cells.Item(148, 7).Value = "Melón"
$ws.Cells.Item(148, 8).Value = "Tuna"
$ws.Cells.Item(148, 9).Value = "Extra"
$ws.Cells.Item(148, 10).Value = 2000
$ws.Cells.Item(148, 11).Value = 800
$ws.Cells.Item(148, 12).Value = 900
$ws.Cells.Item(148, 13).Value = 850
$ws.Cells.Item(148, 14).Value = "$/unidad"
$ws.Cells.Item(148, 15).Value = "Región del Maule"
$ws.Cells.Item(148, 16).Value = 850
$ws.Cells.Item(148, 17).Value = 1
$ws.Cells.Item(148, 18).Value = "Hortaliza"

# --- Row 149: Melón / Tuna / Primera ---
$ws.Cells.Item(149, 1).Value = 7
$ws.Cells.Item(149, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(149, 3).Value = "Ñuble"
$ws.Cells.Item(149, 4).Value = 44568
$ws.Cells.Item(149, 5).Value = 16
$ws.Cells.Item(149, 6).Value = 100112027
$ws.Cells.Item(149, 7).Value = "Melón"
$ws.Cells.Item(149, 8).Value = "Tuna"
$ws.Cells.Item(149, 9).Value = "Primera"
$ws.Cells.Item(149, 10).Value = 4000
$ws.Cells.Item(149, 11).Value = 600
$ws.Cells.Item(149, 12).Value = 700
$ws.Cells.Item(149, 13).Value = 650
$ws.Cells.Item(149, 14).Value = "$/unidad"
$ws.Cells.Item(149, 15).Value = "Región del Maule"
$ws.Cells.Item(149, 16).Value = 650
$ws.Cells.Item(149, 17).Value = 1
$ws.Cells.Item(149, 18).Value = "Hortaliza"
